$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting of the last existing data row (877) down to the new rows (878-920)
$ws.Range("A877:B877").Copy()
$ws.Range("A878:B920").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$aValues = @(876,877,878,879,880,881,882,883,884,885,886,887,888,889,890,891,892,893,894,895,896,897,898,899,900,901,902,903,904,905,906,907,908,909,910,911,912,913,914,915,916,917,918)
$bValues = @(0.30726,0.30726,0.32464,0.31037,0.31595,0.29919,0.30478,0.30168,0.29919,0.2905,0.2905,0.2905,0.28181,0.2396,0.25388,0.23712,0.2396,0.23402,0.22222,0.21664,0.21167,0.20546,0.20608,0.17194,0.18063,0.20732,0.19429,0.19429,0.21788,0.19491,0.2067,0.20608,0.21167,0.21415,0.22222,0.21353,0.2185,0.20857,0.19181,0.21664,0.22843,0.26505,0.27933)

$startRow = 878
for ($i = 0; $i -lt $aValues.Length; $i++) {
    $r = $startRow + $i
    $ws.Cells.Item($r, 1).Value = $aValues[$i]
    $ws.Cells.Item($r, 2).Value = $bValues[$i]
}
